$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bets")

# Insert a new row 28 right after the last data row (27). Inserting (rather
# than just writing into the blank row below the used range) makes Excel
# carry the formatting of the row above down into the new row, same as
# typing a new entry directly under the existing table.
$ws.Rows.Item(28).Insert()

# --- New bet entry (row 28) ---
$ws.Cells.Item(28, 1).Value = 27
$ws.Cells.Item(28, 2).Value = 45221
$ws.Cells.Item(28, 3).Value = 1
$ws.Cells.Item(28, 4).Formula = "=F27"
$ws.Cells.Item(28, 5).Value = 2010
$ws.Cells.Item(28, 6).Formula = "=D28+E28"
$ws.Cells.Item(28, 7).Value = "ESPORTS"
$ws.Cells.Item(28, 8).Value = "WORLDS 2023"
$ws.Cells.Item(28, 9).Value = "BLG"
$ws.Cells.Item(28, 10).Value = "PRIMER INHIBIDOR"
$ws.Cells.Item(28, 11).Value = 1
$ws.Cells.Item(28, 12).Value = 0
$ws.Cells.Item(28, 13).Formula = "=ROUND((F28/`$D`$2-1)*100, 3)"

# --- Update the view state to match where the user ended up scrolled/selected ---
$ws.Range("E29").Select()
$excel.ActiveWindow.ScrollRow = 12
$excel.ActiveWindow.ScrollColumn = 6
